$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.493.76'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.819.96'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.21'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5092'
$ws.Range("E7").Value = '  -6.75%  '
$ws.Range("E8").Value = '  -2.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08163'
$ws.Range("E9").Value = '  +5.98%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.66'
$ws.Range("E10").Value = '  -0.54%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.108'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.340'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.10'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.526'
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("D16").Value = '1.815.16'
$ws.Range("E16").Value = '  -0.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001135'
$ws.Range("E17").Value = '  +4.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.37'
$ws.Range("E18").Value = '  +2.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06654'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.82'
$ws.Range("E20").Value = '  +0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.096'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").Value = '28.522.87'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.42'
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.268'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.33'
$ws.Range("E26").Value = '  +2.80%  '
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("D28").Value = '2.025.94'
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("E29").Value = '  -2.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.82'
$ws.Range("E30").Value = '  +1.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.112'
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1096'
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.771'
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07064'
$ws.Range("E35").Value = '  -4.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2224'
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02349'
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.231'
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.832'
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6311'
$ws.Range("E40").Value = '  +0.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.28'
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.180'
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.402'
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.52'
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.736'
$ws.Range("E46").Value = '  +1.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5921'
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.97'
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.985'
$ws.Range("E49").Value = '  -1.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.182'
$ws.Range("E50").Value = '  -1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06891'
$ws.Range("E51").Value = '  +0.33%  '
